# This script applies two changes to the workbook:
# 1. Updates the F-column "time_taken" timestamps on the "data" sheet
#    to reflect a later re-run of the panel query.
# 2. Adds a new "metadata" worksheet (placed after "data") that records
#    provenance information about the panel query.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$timestamps = @{
    2 = "2021-10-05 14:22:56.631108"
    3 = "2021-10-05 14:22:56.631115"
    4 = "2021-10-05 14:22:56.631119"
    5 = "2021-10-05 14:22:56.631121"
    6 = "2021-10-05 14:22:56.631124"
    7 = "2021-10-05 14:22:56.631126"
    8 = "2021-10-05 14:22:56.631129"
    9 = "2021-10-05 14:22:56.631131"
    10 = "2021-10-05 14:22:56.631134"
    11 = "2021-10-05 14:22:56.631137"
    12 = "2021-10-05 14:22:56.631139"
    13 = "2021-10-05 14:22:56.631141"
    14 = "2021-10-05 14:22:56.631144"
    15 = "2021-10-05 14:22:56.631146"
    16 = "2021-10-05 14:22:56.631149"
    17 = "2021-10-05 14:22:56.631151"
    18 = "2021-10-05 14:22:56.631154"
    19 = "2021-10-05 14:22:56.631156"
    20 = "2021-10-05 14:22:56.631159"
    21 = "2021-10-05 14:22:56.631161"
    22 = "2021-10-05 14:22:56.631164"
    23 = "2021-10-05 14:22:56.631166"
    24 = "2021-10-05 14:22:56.631168"
    25 = "2021-10-05 14:22:56.631171"
    26 = "2021-10-05 14:22:56.631174"
    27 = "2021-10-05 14:22:56.631176"
    28 = "2021-10-05 14:22:56.631179"
    29 = "2021-10-05 14:22:56.631181"
    30 = "2021-10-05 14:22:56.631183"
    31 = "2021-10-05 14:22:56.631186"
    32 = "2021-10-05 14:22:56.631188"
    33 = "2021-10-05 14:22:56.631191"
    34 = "2021-10-05 14:22:56.631194"
    35 = "2021-10-05 14:22:56.631196"
    36 = "2021-10-05 14:22:56.631199"
    37 = "2021-10-05 14:22:56.631201"
    38 = "2021-10-05 14:22:56.631204"
    39 = "2021-10-05 14:22:56.631206"
    40 = "2021-10-05 14:22:56.631208"
    41 = "2021-10-05 14:22:56.631211"
    42 = "2021-10-05 14:22:56.631214"
    43 = "2021-10-05 14:22:56.631216"
    44 = "2021-10-05 14:22:56.631219"
    45 = "2021-10-05 14:22:56.631221"
    46 = "2021-10-05 14:22:56.631223"
    47 = "2021-10-05 14:22:56.631226"
    48 = "2021-10-05 14:22:56.631228"
    49 = "2021-10-05 14:22:56.631231"
    50 = "2021-10-05 14:22:56.631233"
    51 = "2021-10-05 14:22:56.631235"
    52 = "2021-10-05 14:22:56.631238"
    53 = "2021-10-05 14:22:56.631240"
    54 = "2021-10-05 14:22:56.631243"
    55 = "2021-10-05 14:22:56.631245"
    56 = "2021-10-05 14:22:56.631248"
    57 = "2021-10-05 14:22:56.631250"
    58 = "2021-10-05 14:22:56.631253"
    59 = "2021-10-05 14:22:56.631255"
    60 = "2021-10-05 14:22:56.631257"
    61 = "2021-10-05 14:22:56.631260"
    62 = "2021-10-05 14:22:56.631262"
    63 = "2021-10-05 14:22:56.631265"
    64 = "2021-10-05 14:22:56.631267"
    65 = "2021-10-05 14:22:56.631269"
    66 = "2021-10-05 14:22:56.631273"
    67 = "2021-10-05 14:22:56.631276"
    68 = "2021-10-05 14:22:56.631278"
    69 = "2021-10-05 14:22:56.631281"
    70 = "2021-10-05 14:22:56.631283"
    71 = "2021-10-05 14:22:56.631285"
    72 = "2021-10-05 14:22:56.631288"
    73 = "2021-10-05 14:22:56.631290"
    74 = "2021-10-05 14:22:56.631293"
    75 = "2021-10-05 14:22:56.631295"
    76 = "2021-10-05 14:22:56.631298"
    77 = "2021-10-05 14:22:56.631300"
    78 = "2021-10-05 14:22:56.631305"
    79 = "2021-10-05 14:22:56.631308"
    80 = "2021-10-05 14:22:56.631310"
    81 = "2021-10-05 14:22:56.631313"
    82 = "2021-10-05 14:22:56.631315"
    83 = "2021-10-05 14:22:56.631317"
    84 = "2021-10-05 14:22:56.631320"
    85 = "2021-10-05 14:22:56.631322"
    86 = "2021-10-05 14:22:56.631325"
    87 = "2021-10-05 14:22:56.631327"
    88 = "2021-10-05 14:22:56.631330"
    89 = "2021-10-05 14:22:56.631332"
    90 = "2021-10-05 14:22:56.631334"
    91 = "2021-10-05 14:22:56.631337"
    92 = "2021-10-05 14:22:56.631339"
    93 = "2021-10-05 14:22:56.631342"
    94 = "2021-10-05 14:22:56.631345"
    95 = "2021-10-05 14:22:56.631348"
    96 = "2021-10-05 14:22:56.631351"
    97 = "2021-10-05 14:22:56.631353"
    98 = "2021-10-05 14:22:56.631356"
    99 = "2021-10-05 14:22:56.631358"
    100 = "2021-10-05 14:22:56.631361"
    101 = "2021-10-05 14:22:56.631363"
    102 = "2021-10-05 14:22:56.631365"
    103 = "2021-10-05 14:22:56.631368"
    104 = "2021-10-05 14:22:56.631370"
    105 = "2021-10-05 14:22:56.631373"
    106 = "2021-10-05 14:22:56.631375"
    107 = "2021-10-05 14:22:56.631378"
    108 = "2021-10-05 14:22:56.631380"
    109 = "2021-10-05 14:22:56.631382"
    110 = "2021-10-05 14:22:56.631387"
    111 = "2021-10-05 14:22:56.631390"
    112 = "2021-10-05 14:22:56.631392"
    113 = "2021-10-05 14:22:56.631395"
    114 = "2021-10-05 14:22:56.631397"
    115 = "2021-10-05 14:22:56.631399"
    116 = "2021-10-05 14:22:56.631402"
    117 = "2021-10-05 14:22:56.631404"
    118 = "2021-10-05 14:22:56.631407"
    119 = "2021-10-05 14:22:56.631409"
    120 = "2021-10-05 14:22:56.631412"
    121 = "2021-10-05 14:22:56.631414"
    122 = "2021-10-05 14:22:56.631416"
    123 = "2021-10-05 14:22:56.631419"
    124 = "2021-10-05 14:22:56.631421"
    125 = "2021-10-05 14:22:56.631424"
    126 = "2021-10-05 14:22:56.631426"
    127 = "2021-10-05 14:22:56.631428"
    128 = "2021-10-05 14:22:56.631431"
    129 = "2021-10-05 14:22:56.631433"
    130 = "2021-10-05 14:22:56.631438"
    131 = "2021-10-05 14:22:56.631441"
    132 = "2021-10-05 14:22:56.631443"
    133 = "2021-10-05 14:22:56.631446"
    134 = "2021-10-05 14:22:56.631448"
}

foreach ($rowNum in $timestamps.Keys) {
    $ws1.Cells.Item($rowNum, 6).Value = $timestamps[$rowNum]
}

# --- Add the "metadata" worksheet ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "metadata"

# Reuse the header style (bold, centered, bordered) already used by the
# "data" sheet's header row, and the style used for the index column, so
# no new cell styles are introduced.
$ws1.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Thoracic dystrophies"
$newSheet.Range("C2").Value = 122

# data_version must stay a text value ("1.12"), not be coerced to a number.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.12"

$newSheet.Range("E2").Value = "2021-01-29T11:22:46.287036Z"
$newSheet.Range("F2").Value = "2021-10-05 14:22:56.627911"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/122/?format=json"

# Keep "data" as the active/selected sheet, as it was before this edit.
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null

Write-Host "metadata sheet added and timestamps refreshed"
